$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025"
# ---------------------------------------------------------------------
$datePara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "September 19, 2025") {
        $datePara = $p
        break
    }
}
if ($datePara -ne $null) {
    $datePara.Range.Text = "September 21, 2025"
}

# ---------------------------------------------------------------------
# 2. Split the sender's return-address line
#       "919 Story Road, San Jose CA 95122"
#    into two separate paragraphs:
#       "919 Story Road"
#       "San Jose, CA 95122"
#    (The same text also appears later as the "PROPERTY ADDRESS:" value;
#    that occurrence is left untouched, so only the FIRST paragraph with
#    this exact text - the sender's return address, right after
#    "nan The Vu Family Trust" - is split.)
# ---------------------------------------------------------------------
$addrPara = $null
$addrIndex = -1
$idx = 0
foreach ($p in $d.Paragraphs) {
    $idx = $idx + 1
    if ($p.Range.Text.TrimEnd() -eq "919 Story Road, San Jose CA 95122") {
        $addrPara = $p
        $addrIndex = $idx
        break
    }
}

if ($addrPara -ne $null) {
    $addrPara.Range.Text = "919 Story Road"

    # Insert a brand-new paragraph right after it (ahead of the pre-existing
    # blank paragraph that used to follow the address), matching formatting.
    $followingPara = $d.Paragraphs($addrIndex + 1)
    $insertPoint = $d.Range($followingPara.Range.Start, $followingPara.Range.Start)
    $insertPoint.InsertParagraphBefore()

    $cityPara = $d.Paragraphs($addrIndex + 1)
    # Assign with a trailing space first so the run's <w:t> is recognised as
    # needing xml:space="preserve" (matching the canonical OOXML produced by
    # Word for this text), then trim the space back off.
    $cityPara.Range.Text = "San Jose, CA 95122 "
    $cityPara = $d.Paragraphs($addrIndex + 1)
    $cityPara.Range.Text = "San Jose, CA 95122"
}

# ---------------------------------------------------------------------
# 3. Remove the now-unneeded blank "No Spacing" paragraph that immediately
#    followed "... Board of Directors" in the closing signature block.
# ---------------------------------------------------------------------
$paragraphs = @($d.Paragraphs)
for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    if ($paragraphs[$i].Range.Text -match "Board of Directors") {
        $blankPara = $paragraphs[$i + 1]
        if ($blankPara.Range.Text.Trim() -eq "" -and $blankPara.Style.NameLocal -eq "No Spacing") {
            $blankPara.Range.Delete()
        }
        break
    }
}
